$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Determine the last used row in column D (the "enddate" column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

$rng = $ws.Range("D2:D$lastRow")

# Force text formatting before assignment so Excel does not auto-convert
# the "2025-03-11" string into a date serial number (the source cells are
# plain text, not dates). Restore the original "General" formatting
# afterwards so only the cell contents change.
$rng.NumberFormat = "@"
$rng.Value = "2025-03-11"
$rng.ClearFormats()
